$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "entree" -> "entrance" typo.
# Set the Housing_entree -> Housing_entrance id cells first so the new
# shared string "Housing_entrance" is registered before "1. Housing entrance".
$ws.Range("H2").Value = "Housing_entrance"
$ws.Range("K2").Value = "Housing_entrance"
$ws.Range("A2").Value = "1. Housing entrance"

# Row 14 (Atrium) had a leftover cell format applied (fill) that is no
# longer needed - clear it back to the default style.
$ws.Range("A14:K14").ClearFormats()

# Update the active selection to match the saved view state.
$ws.Range("H17").Select()
